$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "id_user" column (A) was removed from the export; deleting the
# entire column shifts B:I left into A:H, which also drops the old
# A2 value ("[boleh kosong]") while preserving the style that used to
# live on C2 (now B2).
$ws.Range("A1").EntireColumn.Delete()

# Column A was resized afterwards (closest width reachable through the
# column-width object model to the authored 19.85546875 char width).
$ws.Columns.Item(1).ColumnWidth = 19

# Selection/view: the sheet now shows column A at the left edge (no more
# scrolled topLeftCell), with the whole of column A selected.
$ws.Range("A1:A1048576").Select() | Out-Null
